# edit.ps1 - applies Review_142 -> Review_141 (DoLa) content update
$d = $word.ActiveDocument

# --- 1) Heading paragraph: title text + huggingface URL ---
$d.Content.Find.Execute('Review 142: LARGE LANGUAGE MODELS AS OPTIMIZERS, 09.09.2023', $true, $false, $false, $false, $false, $true, 1, $false, 'Review 141: [Short] DOLA: DECODING BY CONTRASTING LAYERS IMPROVES FACTUALITY IN LARGE LANGUAGE MODELS, 08.09.2023', 2) | Out-Null
$d.Content.Find.Execute('https://huggingface.co/papers/2309.03409', $true, $false, $false, $false, $false, $true, 1, $false, 'https://huggingface.co/papers/2309.03883', 2) | Out-Null

# --- 2) Bold "Paper:" line with arxiv URL ---
$d.Content.Find.Execute('Paper: https://arxiv.org/abs/2311.15249v1', $true, $false, $false, $false, $false, $true, 1, $false, 'Paper: https://arxiv.org/abs/2309.03883v2', 2) | Out-Null

# --- 3) Replace the body paragraph (5th paragraph) text, keeping its 3 line-break runs ---
$body = $d.Paragraphs(5)
$parts = @(
  'אנחנו משתמשים במודלי שפה למשימות רבות אבל האם לסמוך על פלטיהם? עם כל הצער עדיין לא. יש תופעה הנקראת הזיות (hallucinations) של מודלי שפה כאשר מודלי שפה מדברים שטויות. זו בעיה מאוד רצינית בטח אם אתם רוצים לשים מודל שפה בפרודקשן. ',
  'אז היום ב- #shorthebrewpapereviews אנחנו סוקרים מאמר המציע פתרון לבעיה החמורה הזו. המחברים מבססים את שיטתם על התכונה המעניינת של מודלי שפה שניתן לצפותה כאשר מוסיפים שכבת סופטמקס המחשבת את התפלגות הטוקנים אחרי כל בלוק הטרנספורמר. ',
  'למעשה יש כאן שתי תופעות נפרדות. עבור טוקנים קלים (יחסית) לניחוש (הנובעים מכללי הדקדוק למשל) אז התפלגות טוקנים משכבות האמצע בערך לא משתנה ושווה להתפלגות הסופית של הטוקנים. בטוקנים היותר קשים ההתפלגות משתנה משמעתית כמעט עד השכבה האחרונה – כלומר בשכבה לפני האחרונה התפלגות הטוקנים עשויה להיות שונה מאוד מהשכבה הסופית. ',
  'אחד ההסברים לכך הוא ״שואב את ״הידע העובדתי מהזיכרון שלו״ (סליחה על נפנופי הידיים אבל ככה כתוב במאמר). המחברים מציעים לנרמל את הסופטמקס הסופי עם הסופטמקס של השכבה בעלת שוני הגבוה ביותר (מבחינת התפלגות הטוקנים). כלומר הסתברות התוקן פרופורציאונלית להשתנות המקסימלית של ההסתברות של הטוקן הזה (= ״כמות הלמידה״?). ',
  'עושים זאת רק לטוקנים בעלי הסתברות גבוהה מספיק בשכבת סופטמקס האחרונה – השכבות עם שוני מקסימלי נבחרות על סט ולידציה. השוני נמדד במונחי Jensen-Shannon Divergence או JSD (גרסה סימטרית של KL) בין התפלגויות הטוקנים. מכיוון שנרמול זה עלול לדפוק את הדקדוק מוסיפים ״קנס על חזרתיות״ (שלא יפלוט את אותו הטקסט יותר מפעם אחת).'
)
$lineBreak = [char]11
$newBodyText = $parts[0] + $lineBreak + $lineBreak + $parts[1] + $lineBreak + $lineBreak + $parts[2] + $lineBreak + $lineBreak + $parts[3] + $lineBreak + $lineBreak + $parts[4]
$bodyRange = $body.Range
$bodyRange.End = $bodyRange.End - 1
$bodyRange.Text = $newBodyText

# --- 4) Delete the next three paragraphs (now paragraphs 6, 6, 6 after each delete) ---
$d.Paragraphs(6).Range.Delete() | Out-Null
$d.Paragraphs(6).Range.Delete() | Out-Null
$d.Paragraphs(6).Range.Delete() | Out-Null

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
